$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so that values like
# "55.40" or "0.719" are not auto-converted to numbers by Excel,
# matching the original inlineStr / shared-string cell content.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '65.151.23'
$ws.Range('E2').Value = '  +1.99%  '
$ws.Range('D3').Value = '3.161.95'
$ws.Range('E3').Value = '  +3.25%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '576.87'
$ws.Range('E5').Value = '  +3.36%  '
$ws.Range('D6').Value = '150.30'
$ws.Range('E6').Value = '  +5.46%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '3.158.32'
$ws.Range('E8').Value = '  +3.22%  '
$ws.Range('E9').Value = '  +2.25%  '
$ws.Range('D10').Value = '0.160'
$ws.Range('E10').Value = '  +4.73%  '
$ws.Range('E11').Value = '  -0.46%  '
$ws.Range('D12').Value = '0.501'
$ws.Range('E12').Value = '  +4.54%  '
$ws.Range('E13').Value = '  +14.70%  '
$ws.Range('D14').Value = '37.29'
$ws.Range('E14').Value = '  +5.97%  '
$ws.Range('D15').Value = '3.676.27'
$ws.Range('E15').Value = '  +3.15%  '
$ws.Range('D16').Value = '65.138.11'
$ws.Range('E16').Value = '  +1.93%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.156.55'
$ws.Range('E17').Value = '  +3.06%  '
$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D18').Value = '7.13'
$ws.Range('E18').Value = '  +5.25%  '
$ws.Range('E19').Value = '  +1.06%  '
$ws.Range('D20').Value = '510.05'
$ws.Range('E20').Value = '  +4.08%  '
$ws.Range('D21').Value = '14.84'
$ws.Range('E21').Value = '  +3.82%  '
$ws.Range('B22').Value = 'Polygon'
$ws.Range('C22').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D22').Value = '0.719'
$ws.Range('E22').Value = '  +4.89%  '
$ws.Range('B23').Value = 'InternetComputer(DFINITY)'
$ws.Range('C23').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D23').Value = '15.33'
$ws.Range('D24').Value = '7.76'
$ws.Range('E24').Value = '  +3.27%  '
$ws.Range('D25').Value = '84.55'
$ws.Range('E25').Value = '  +2.10%  '
$ws.Range('E26').Value = '  +0.09%  '
$ws.Range('B27').Value = 'PancakeSwap'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D27').Value = '2.91'
$ws.Range('E27').Value = '  +3.87%  '
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').Value = '8.92'
$ws.Range('E28').Value = '  +9.87%  '
$ws.Range('E29').Value = '  +6.73%  '
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').Value = '27.78'
$ws.Range('E30').Value = '  +4.91%  '
$ws.Range('B31').Value = 'Stacks'
$ws.Range('C31').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D31').Value = '2.78'
$ws.Range('E31').Value = '  +11.55%  '
$ws.Range('D32').Value = '0.999'
$ws.Range('E32').Value = '  -0.06%  '
$ws.Range('E33').Value = '  +3.35%  '
$ws.Range('D34').Value = '6.28'
$ws.Range('E34').Value = '  +10.64%  '
$ws.Range('D35').Value = '6.55'
$ws.Range('E35').Value = '  +5.51%  '
$ws.Range('D36').Value = '55.40'
$ws.Range('E36').Value = '  +0.13%  '
$ws.Range('D37').Value = '0.0903'
$ws.Range('E37').Value = '  +10.73%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '0.0427'
$ws.Range('E38').Value = '  +3.95%  '
$ws.Range('B39').Value = 'Bittensor'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D39').Value = '468.45'
$ws.Range('E39').Value = '  +5.67%  '
$ws.Range('D40').Value = '3.04'
$ws.Range('E40').Value = '  +9.02%  '
$ws.Range('D41').Value = '8.68'
$ws.Range('E41').Value = '  +4.25%  '
$ws.Range('D42').Value = '3.057.32'
$ws.Range('E42').Value = '  +0.90%  '
$ws.Range('E43').Value = '  +1.15%  '
$ws.Range('E44').Value = '  +10.14%  '
$ws.Range('E45').Value = '  +3.80%  '
$ws.Range('D46').Value = '28.83'
$ws.Range('E46').Value = '  +3.77%  '
$ws.Range('D47').Value = '0.0₃0593'
$ws.Range('E47').Value = '  +15.08%  '
$ws.Range('E48').Value = '  -0.03%  '
$ws.Range('D49').Value = '0.115'
$ws.Range('E49').Value = '  +0.89%  '
$ws.Range('D50').Value = '2.26'
$ws.Range('E50').Value = '  +6.46%  '
$ws.Range('D51').Value = '120.04'
$ws.Range('E51').Value = '  +1.49%  '

# Restore the default (unstyled) cell style for column D data cells,
# matching the original workbook where these cells carried no style index.
$ws.Range("D2:D51").Style = "Normal"
